# chore: adapt column header formatting to respective input file names
#
# Renames the two header-row suffixes ("_old" -> "_FV2404", "_new" -> "_FV2410"),
# turns the data range A1:U57 into a native Excel Table (ListObject) bound to
# those (now renamed) headers, and freezes the header row (row 1) in the
# worksheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.Range("A1:U57")
$headerRow = $ws.Range("A1:U1")
$lastCol = $headerRow.Columns.Count

# --- Rename the header labels --------------------------------------------
# Columns A:J carry the "<field>_old" headers -> "<field>_FV2404"
# Column  K   carries "diff" and stays untouched
# Columns L:U carry the "<field>_new" headers -> "<field>_FV2410"
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $text = [string]$cell.Value2

    if ($text -like "*_old") {
        $cell.Value2 = ($text.Substring(0, $text.Length - 4) + "_FV2404")
    } elseif ($text -like "*_new") {
        $cell.Value2 = ($text.Substring(0, $text.Length - 4) + "_FV2410")
    }
}

# --- Convert the range into a proper Table (ListObject) -------------------
$table = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $usedRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$table.Name = "Table1"

# --- Freeze the header row (row 1) ----------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
